# Adds the "2022-Q3" sheet (new quarterly fund-holdings snapshot) and
# updates the "总计" (totals) sheet with the corresponding summary row,
# shifting the existing quarters down by one position.
#
#   总计        : unchanged layout, new row inserted for 2022-Q3, others shift down
#   2022-Q3     : brand-new sheet, placed between 总计 and 2022-Q2
#   2022-Q2 / 2022-Q1 / 2021-Q4 : untouched, just shift one tab to the right

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Build the new "2022-Q3" sheet.
#    Duplicate the existing "2022-Q2" sheet (placing the copy directly
#    before it) so the new sheet inherits the exact same look & feel:
#    bold/bordered header row, bold/bordered index column, plain data
#    cells. Then overwrite its contents with the 2022-Q3 figures.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$ws = $wb.Worksheets.Item(2)
$ws.Name = "2022-Q3"

# The 2022-Q3 snapshot has 9 fund rows (rows 2-10) vs. 2022-Q2's 6 (rows
# 2-7), so extend the row-7 formatting (index-column border/bold, plain
# data columns) down through row 10 before writing the extra rows.
$ws.Range("A7:H7").Copy()
$ws.Range("A8:H10").PasteSpecial(-4122)

# The fund-detail columns (B:G) hold numeric-looking values ("540006",
# "14.62", ...) that must stay TEXT, matching the source convention -
# force Text format before writing so Excel doesn't coerce them to
# numbers.
$ws.Range("B2:G10").NumberFormat = "@"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "540006"
$ws.Range("C2").Value = "汇丰晋信大盘股票A"
$ws.Range("D2").Value = "14.62"
$ws.Range("E2").Value = "94.47"
$ws.Range("F2").Value = "2.90"
$ws.Range("G2").Value = "0.4240"
$ws.Range("H2").Value = 2
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "000965"
$ws.Range("C3").Value = "汇丰晋信新动力混合"
$ws.Range("D3").Value = "0.89"
$ws.Range("E3").Value = "93.15"
$ws.Range("F3").Value = "5.58"
$ws.Range("G3").Value = "0.0497"
$ws.Range("H3").Value = 2
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "540004"
$ws.Range("C4").Value = "汇丰晋信2026周期混合"
$ws.Range("D4").Value = "1.10"
$ws.Range("E4").Value = "33.59"
$ws.Range("F4").Value = "2.90"
$ws.Range("G4").Value = "0.0319"
$ws.Range("H4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "011997"
$ws.Range("C5").Value = "景顺长城安盈回报一年持有期混合A"
$ws.Range("D5").Value = "1.69"
$ws.Range("E5").Value = "26.07"
$ws.Range("F5").Value = "1.34"
$ws.Range("G5").Value = "0.0226"
$ws.Range("H5").Value = 8
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "004557"
$ws.Range("C6").Value = "北信瑞丰鼎丰灵活配置混合"
$ws.Range("D6").Value = "0.32"
$ws.Range("E6").Value = "64.59"
$ws.Range("F6").Value = "3.27"
$ws.Range("G6").Value = "0.0105"
$ws.Range("H6").Value = 9
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "519099"
$ws.Range("C7").Value = "新华灵活主题混合"
$ws.Range("D7").Value = "0.15"
$ws.Range("E7").Value = "81.48"
$ws.Range("F7").Value = "2.14"
$ws.Range("G7").Value = "0.0032"
$ws.Range("H7").Value = 9
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "005966"
$ws.Range("C8").Value = "安信中证500指数增强C"
$ws.Range("D8").Value = "0.16"
$ws.Range("E8").Value = "92.50"
$ws.Range("F8").Value = "1.18"
$ws.Range("G8").Value = "0.0019"
$ws.Range("H8").Value = 4
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "005965"
$ws.Range("C9").Value = "安信中证500指数增强A"
$ws.Range("D9").Value = "0.10"
$ws.Range("E9").Value = "92.50"
$ws.Range("F9").Value = "1.18"
$ws.Range("G9").Value = "0.0012"
$ws.Range("H9").Value = 4
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "011998"
$ws.Range("C10").Value = "景顺长城安盈回报一年持有期混合C"
$ws.Range("D10").Value = "0.08"
$ws.Range("E10").Value = "26.07"
$ws.Range("F10").Value = "1.34"
$ws.Range("G10").Value = "0.0011"
$ws.Range("H10").Value = 8

# ---------------------------------------------------------------------
# 2. Update "总计" with the new 2022-Q3 summary row; existing rows shift
#    down by one (index column re-numbered 0..3).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Extend the bold/bordered index-column formatting to the new row 5.
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 9
$total.Range("D2").Value = 0.55

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 6
$total.Range("D3").Value = 0.52

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 50
$total.Range("D4").Value = 5.97

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 18
$total.Range("D5").Value = 8.07
